$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 144, shifting the existing rows 144:204 down to 145:205
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with the new record
$ws.Range("A144").Value = 10
$ws.Range("B144").Value = "Vega Modelo de Temuco"
$ws.Range("C144").Value = "La Araucanía"
$ws.Range("D144").Value = 44510
$ws.Range("E144").Value = 9
$ws.Range("F144").Value = 100112017
$ws.Range("G144").Value = "Apio"
$ws.Range("H144").Value = "Americana (o)"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 40
$ws.Range("K144").Value = 9000
$ws.Range("L144").Value = 9000
$ws.Range("M144").Value = 9000
$ws.Range("N144").Value = "$/docena de matas"
$ws.Range("O144").Value = "Región Metropolitana"
$ws.Range("P144").Value = 1500
$ws.Range("Q144").Value = 6
$ws.Range("R144").Value = "Hortaliza"
